$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 724
$ws1.Range("F6").Value = 2932
$ws1.Range("F8").Value = 1943
$ws1.Range("F9").Value = 318
$ws1.Range("F10").Value = 291
$ws1.Range("F11").Value = 792
$ws1.Range("F12").Value = 936
$ws1.Range("F14").Value = 400
$ws1.Range("F15").Value = 1131
$ws1.Range("F17").Value = 61
$ws1.Range("F19").Value = 7032
$ws1.Range("F20").Value = 266
$ws1.Range("F21").Value = 1733
$ws1.Range("F22").Value = 181
$ws1.Range("F25").Value = 359
$ws1.Range("F26").Value = 287
$ws1.Range("F27").Value = 75
$ws1.Range("F29").Value = 932
$ws1.Range("F30").Value = 63
$ws1.Range("F31").Value = 111
$ws1.Range("F34").Value = 1918
$ws1.Range("F35").Value = 166
$ws1.Range("F38").Value = 244
$ws1.Range("F39").Value = 30
$ws1.Range("F41").Value = 256
$ws1.Range("F43").Value = 189

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 724
$ws4.Range("F9").Value = 2932
$ws4.Range("F11").Value = 1943
$ws4.Range("F12").Value = 318
$ws4.Range("F13").Value = 291
$ws4.Range("F14").Value = 792
$ws4.Range("F16").Value = 936
$ws4.Range("F18").Value = 400
$ws4.Range("F19").Value = 1131
$ws4.Range("F20").Value = 61
$ws4.Range("F22").Value = 7032
$ws4.Range("F23").Value = 266
$ws4.Range("F24").Value = 1733
$ws4.Range("F26").Value = 181
$ws4.Range("F29").Value = 359
$ws4.Range("F30").Value = 287
$ws4.Range("F31").Value = 75
$ws4.Range("F33").Value = 932
$ws4.Range("F34").Value = 63
$ws4.Range("F35").Value = 111
$ws4.Range("F37").Value = 1918
$ws4.Range("F38").Value = 166
$ws4.Range("F41").Value = 244
$ws4.Range("F42").Value = 30
$ws4.Range("F44").Value = 256
$ws4.Range("F49").Value = 189
